$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at 653, pushing the existing rows 653-658 down to 658-663
$ws.Rows("653:657").Insert()

# Common columns shared by every row in this dataset
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$tipo = "Fruta"
$productoId = 100101
$producto = "Berries"
$categoriaId = 100101007
$categoria = "Kiwi"
$unidad = "`$/bins (450 kilos)"
$origen = "Región de O'Higgins"
$kgUnidad = 450

function Set-KiwiRow {
    param(
        [int]$row,
        [double]$fecha,
        [string]$variedad,
        [string]$calidad,
        [double]$volumen,
        [double]$precioMin,
        [double]$precioMax,
        [double]$precioProm,
        [double]$precioKg
    )

    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-KiwiRow 653 44656 "Gold"    "Primera"                8  120000 120000 120000 267
Set-KiwiRow 654 44656 "Gold"    "Segunda"                12 100000 100000 100000 222
Set-KiwiRow 655 44656 "Hayward" "Extra (doble especial)" 15 430000 430000 430000 956
Set-KiwiRow 656 44656 "Hayward" "Primera"                36 260000 280000 270000 600
Set-KiwiRow 657 44656 "Hayward" "Segunda"                22 230000 230000 230000 511
